$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1518.7894
$ws.Range("I28").Value = 900.4167
$ws.Range("J28").Value = 2578.8572
$ws.Range("K28").Value = 900.4167
$ws.Range("L28").Value = 2578.8572
$ws.Range("M28").Value = -415.4167
$ws.Range("N28").Value = -3548.8572
$ws.Range("H69").Value = 3417.4285
$ws.Range("I69").Value = 3107.6
$ws.Range("J69").Value = 3589.5557
$ws.Range("K69").Value = 9322.799999999999
$ws.Range("L69").Value = 10768.6671
$ws.Range("M69").Value = -8448.799999999999
$ws.Range("N69").Value = -12516.6671
$ws.Range("H72").Value = 3417.4285
$ws.Range("I72").Value = 3107.6
$ws.Range("J72").Value = 3589.5557
$ws.Range("K72").Value = 27968.4
$ws.Range("L72").Value = 32306.0013
$ws.Range("M72").Value = -23600.4
$ws.Range("N72").Value = -41042.0013
$ws.Range("H74").Value = 3399.7104
$ws.Range("I74").Value = 3026.8635
$ws.Range("K74").Value = 3026.8635
$ws.Range("M74").Value = -2090.8635
$ws.Range("H77").Value = 3399.7104
$ws.Range("I77").Value = 3026.8635
$ws.Range("K77").Value = 15134.3175
$ws.Range("M77").Value = -10454.3175
$ws.Range("H129").Value = 1914.8948
$ws.Range("J129").Value = 2570
$ws.Range("L129").Value = 7710
$ws.Range("N129").Value = -17710
$ws.Range("H132").Value = 3524229.8
$ws.Range("I132").Value = 3184.772
$ws.Range("J132").Value = 17859912
$ws.Range("K132").Value = 9554.315999999999
$ws.Range("L132").Value = 53579736
$ws.Range("M132").Value = -7024.315999999999
$ws.Range("N132").Value = -53584796
$ws.Range("H135").Value = 10870431
$ws.Range("I135").Value = 828.8095
$ws.Range("J135").Value = 125001250
$ws.Range("K135").Value = 7459.2855
$ws.Range("L135").Value = 1125011250
$ws.Range("M135").Value = -4924.2855
$ws.Range("N135").Value = -1125016320
$ws.Range("H137").Value = 7143793
$ws.Range("I137").Value = 692.5
$ws.Range("J137").Value = 14286893
$ws.Range("K137").Value = 2077.5
$ws.Range("L137").Value = 42860679
$ws.Range("M137").Value = 472.5
$ws.Range("N137").Value = -42865779

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10716.38
$ws.Range("I32").Value = 10337
$ws.Range("J32").Value = 12335.066
$ws.Range("K32").Value = 10337
$ws.Range("L32").Value = 12335.066
$ws.Range("M32").Value = -10050
$ws.Range("N32").Value = -12909.066
$ws.Range("H61").Value = 18520826
$ws.Range("I61").Value = 20002056
$ws.Range("J61").Value = 5450
$ws.Range("K61").Value = 20002056
$ws.Range("L61").Value = 5450
$ws.Range("M61").Value = -20001844
$ws.Range("N61").Value = -5874
$ws.Range("H74").Value = 7814660.5
$ws.Range("I74").Value = 9616843
$ws.Range("J74").Value = 5203.5
$ws.Range("K74").Value = 9616843
$ws.Range("L74").Value = 5203.5
$ws.Range("M74").Value = -9615969
$ws.Range("N74").Value = -6951.5
$ws.Range("H77").Value = 7814660.5
$ws.Range("I77").Value = 9616843
$ws.Range("J77").Value = 5203.5
$ws.Range("K77").Value = 48084215
$ws.Range("L77").Value = 26017.5
$ws.Range("M77").Value = -48079847
$ws.Range("N77").Value = -34753.5
$ws.Range("H136").Value = 18520826
$ws.Range("I136").Value = 20002056
$ws.Range("J136").Value = 5450
$ws.Range("K136").Value = 60006168
$ws.Range("L136").Value = 16350
$ws.Range("M136").Value = -60003618
$ws.Range("N136").Value = -21450

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2602.3022
$ws.Range("I134").Value = 1571.4166
$ws.Range("J134").Value = 7904
$ws.Range("K134").Value = 4714.2498
$ws.Range("L134").Value = 23712
$ws.Range("M134").Value = -2179.2498
$ws.Range("N134").Value = -28782

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12353065
$ws.Range("I31").Value = 8600.947
$ws.Range("J31").Value = 41671170
$ws.Range("K31").Value = 8600.947
$ws.Range("L31").Value = 41671170
$ws.Range("M31").Value = -8305.947
$ws.Range("N31").Value = -41671760
$ws.Range("H34").Value = 12353065
$ws.Range("I34").Value = 8600.947
$ws.Range("J34").Value = 41671170
$ws.Range("K34").Value = 8600.947
$ws.Range("L34").Value = 41671170
$ws.Range("M34").Value = -8398.947
$ws.Range("N34").Value = -41671574
$ws.Range("H99").Value = 1828.5
$ws.Range("I99").Value = 1150
$ws.Range("J99").Value = 2507
$ws.Range("K99").Value = 1150
$ws.Range("L99").Value = 2507
$ws.Range("M99").Value = 348
$ws.Range("N99").Value = -5503
$ws.Range("H107").Value = 471.86957
$ws.Range("I107").Value = 459.6111
$ws.Range("J107").Value = 516
$ws.Range("K107").Value = 459.6111
$ws.Range("L107").Value = 516
$ws.Range("M107").Value = 1460.3889
$ws.Range("N107").Value = -4356
$ws.Range("H120").Value = 39850.332
$ws.Range("J120").Value = 39850.332
$ws.Range("L120").Value = 39850.332
$ws.Range("N120").Value = -47108.332
$ws.Range("H126").Value = 1828.5
$ws.Range("I126").Value = 1150
$ws.Range("J126").Value = 2507
$ws.Range("K126").Value = 3450
$ws.Range("L126").Value = 7521
$ws.Range("M126").Value = -980
$ws.Range("N126").Value = -12461
$ws.Range("H132").Value = 6850776
$ws.Range("I132").Value = 8475861
$ws.Range("J132").Value = 2202.4285
$ws.Range("K132").Value = 25427583
$ws.Range("L132").Value = 6607.2855
$ws.Range("M132").Value = -25425053
$ws.Range("N132").Value = -11667.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 650
$ws.Range("I36").Value = 650
$ws.Range("K36").Value = 1950
$ws.Range("M36").Value = -1781
$ws.Range("H131").Value = 802.27
$ws.Range("I131").Value = 453.92856
$ws.Range("J131").Value = 858.97675
$ws.Range("K131").Value = 1361.78568
$ws.Range("L131").Value = 2576.93025
$ws.Range("M131").Value = 3678.21432
$ws.Range("N131").Value = -12656.93025

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10420048
$ws.Range("I80").Value = 25643764
$ws.Range("J80").Value = 3820
$ws.Range("K80").Value = 25643764
$ws.Range("L80").Value = 3820
$ws.Range("M80").Value = -25642766
$ws.Range("N80").Value = -5816
$ws.Range("H83").Value = 10420048
$ws.Range("I83").Value = 25643764
$ws.Range("J83").Value = 3820
$ws.Range("K83").Value = 128218820
$ws.Range("L83").Value = 19100
$ws.Range("M83").Value = -128213828
$ws.Range("N83").Value = -29084

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4251.827
$ws.Range("I122").Value = 3967.353
$ws.Range("J122").Value = 4789.1665
$ws.Range("K122").Value = 11902.059
$ws.Range("L122").Value = 14367.4995
$ws.Range("M122").Value = -9452.059000000001
$ws.Range("N122").Value = -19267.4995
$ws.Range("H132").Value = 5887766
$ws.Range("I132").Value = 3249.7576
$ws.Range("J132").Value = 26328716
$ws.Range("K132").Value = 9749.272799999999
$ws.Range("L132").Value = 78986148
$ws.Range("M132").Value = -7219.272799999999
$ws.Range("N132").Value = -78991208
$ws.Range("H136").Value = 7465212
$ws.Range("I136").Value = 8475521
$ws.Range("J136").Value = 14188.125
$ws.Range("K136").Value = 25426563
$ws.Range("L136").Value = 42564.375
$ws.Range("M136").Value = -25424013
$ws.Range("N136").Value = -47664.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 612.4483
$ws.Range("I107").Value = 706.26086
$ws.Range("J107").Value = 252.83333
$ws.Range("K107").Value = 2118.78258
$ws.Range("L107").Value = 758.49999
$ws.Range("M107").Value = -198.7825800000001
$ws.Range("N107").Value = -4598.49999
$ws.Range("H113").Value = 1257.375
$ws.Range("I113").Value = 813.41174
$ws.Range("J113").Value = 2335.5715
$ws.Range("K113").Value = 2440.23522
$ws.Range("L113").Value = 7006.7145
$ws.Range("M113").Value = -270.23522
$ws.Range("N113").Value = -11346.7145
$ws.Range("H132").Value = 1982.4222
$ws.Range("I132").Value = 1556.6923
$ws.Range("J132").Value = 4749.6665
$ws.Range("K132").Value = 4670.0769
$ws.Range("L132").Value = 14248.9995
$ws.Range("M132").Value = -2140.0769
$ws.Range("N132").Value = -19308.9995
$ws.Range("H136").Value = 1207.5264
$ws.Range("I136").Value = 1032.3704
$ws.Range("J136").Value = 1637.4546
$ws.Range("K136").Value = 3097.1112
$ws.Range("L136").Value = 4912.3638
$ws.Range("M136").Value = -547.1112000000003
$ws.Range("N136").Value = -10012.3638
